# regen sval data to filter save games
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @(3.286832544864788, 1.655778082260271, 0.1494219747398047, 0.4942365360607697, 1, 5.586269137925634)
    3  = @(3.286832544864788, 1.655778082260271, 0.1494219747398047, 0.4942365360607697, 1, 5.586269137925634)
    4  = @(3.286832544864788, 1.655778082260271, 0.7527432677738641, 0.4942365360607697, 1, 6.189590430959694)
    5  = @(0.01293466051926884, 1.655778082260271, 22.3905356188092, 10.19245300693656, 1, 34.25170136852529)
    6  = @(3.286832544864788, 1.655778082260271, 3.537761648806719, 0.4942365360607697, 1, 8.974608811992548)
    7  = @(0.6606524410359556, 0.04071648406533734, 0.7527432677738641, 0.4942365360607697, 0, 1.948348728935927)
    8  = @(1.455362044514542, 1.655778082260271, 0.7527432677738641, 0.4942365360607697, 1, 4.358119930609447)
    9  = @(0.1190320826869504, 0.306821227259698, 0.7527432677738641, 0.4942365360607697, 1, 1.672833113781282)
    10 = @(3.286832544864788, 1.655778082260271, 0.1494219747398047, 0.4942365360607697, 0, 5.586269137925634)
    11 = @(0.1190320826869504, 1.655778082260271, 0.1494219747398047, 0.4942365360607697, 0, 2.418468675747795)
    12 = @(1.455362044514542, 1.655778082260271, 3.537761648806719, 0.4942365360607697, 1, 7.143138311642302)
    13 = @(0.2917716402565462, 0.306821227259698, 0.7527432677738641, 0.4942365360607697, 1, 1.845572671350878)
    14 = @(0.04271373187048222, 1.655778082260271, 0.1494219747398047, 0.4942365360607697, 1, 2.342150324931327)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $col = 2
    foreach ($v in $vals) {
        $ws.Cells.Item($row, $col).Value = $v
        $col++
    }
}
